# pr-approval-process.pptx — "Fixing event sequence and modify outline style"
#
# 1) Give a solid (1.5pt / 19050 EMU) outline to the pink (E36386) process
#    boxes that previously had a line with no explicit weight (they inherit
#    the theme default, which renders hairline-thin). This affects the
#    7 contributor/approval rectangles.
# 2) Bump the step counter in the last callout oval from "5" to "6" to fix
#    the event sequence numbering.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shapes whose outline needs an explicit 1.5pt (19050 EMU) weight.
$outlineShapeIds = @(14, 32, 51, 69, 70, 71, 83)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)

    if ($outlineShapeIds -contains $shp.Id) {
        $shp.Line.Weight = 1.5
    }

    if ($shp.Id -eq 89) {
        $shp.TextFrame.TextRange.Text = "6"
    }
}
